$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (gitlab_project_namespace) is now the same value for every row.
$ws.Range("B2").Value = "repo-migration"
$ws.Range("B3").Value = "repo-migration"
$ws.Range("B4").Value = "repo-migration"
$ws.Range("B5").Value = "repo-migration"
$ws.Range("B6").Value = "repo-migration"

# Column C (project_to_import) - new project list.
$ws.Range("C2").Value = "almatasks"
$ws.Range("C3").Value = "app-n-pak"
$ws.Range("C4").Value = "casa-build-utils"
$ws.Range("C5").Value = "casa6"
$ws.Range("C6").Value = "casashell"

# Column D (azure_target_namespace) - same target for every row.
$ws.Range("D2").Value = "repo-migartion/git-project"
$ws.Range("D3").Value = "repo-migartion/git-project"
$ws.Range("D4").Value = "repo-migartion/git-project"
$ws.Range("D5").Value = "repo-migartion/git-project"
$ws.Range("D6").Value = "repo-migartion/git-project"

# Column A (sr) - add the two new row numbers.
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5

# Update selection to match the final state (B6 selected).
$ws.Range("B6").Select()
